# adicionado utils, pathlib e status atrasado para nao enviados
#
# For every row whose "Situacao" (column E) is still blank (i.e. the
# technician's form hasn't been sent), mark it as "Atrasado" (late) with a
# distinctive orange fill (matching the existing "Enviado" cells' white
# Arial font, just on an orange background instead of dark green).
#
# Also bumps the two "Data de Envio" (column F) values that moved since the
# last snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column E ("Situacao") that are currently empty and must become
# "Atrasado".
$atrasadoRows = @(
    2,3,4,6,9,11,12,13,14,15,16,17,19,22,23,24,25,26,27,28,29,30,31,32,
    37,38,39,40,41,42,43,45,47,48,49,51,53,54,55,56,58,60,61,62,63,64,67,68
)

foreach ($r in $atrasadoRows) {
    $cell = $ws.Cells.Item($r, 5)   # column E
    $cell.Value = "Atrasado"

    # Same white Arial used by the "Enviado" cells, just on an orange fill
    # instead of dark green, so "Atrasado" reads as its own status colour.
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 11
    $cell.Font.Bold = $false
    $cell.Font.Color = 16777215        # white
    $cell.Interior.Color = 25855       # RGB(255,100,0) -> orange

    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.VerticalAlignment = -4108    # xlCenter
}

# Two technicians' submission dates were updated.
$ws.Cells.Item(18, 6).Value = "01/04/2025"   # F18 (Douglas de Mattia / Alvaro Rodrigues)
$ws.Cells.Item(59, 6).Value = "31/03/2025"   # F59 (Adnan navarro de freitas kassim)

Write-Output "Atrasado status applied to $($atrasadoRows.Count) rows; F18 and F59 dates updated."
